# "takes readings, writes readings"
# - clear out the old placeholder readings on Sheet1
# - add a new "Manual" sheet (placed right after Sheet1) with the
#   Throttle / Thrust / Current header row that manual readings get logged under
# - leave Sheet1's selection parked away from the data, and make "Manual" the
#   active sheet/selection, matching the saved workbook view state

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Wipe the old sample readings (rows 1-4, cols A-C) from Sheet1.
$ws1.Range("A1:C4").ClearContents()

# Insert the new "Manual" sheet directly after Sheet1.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Manual"

# Header row for logging throttle/thrust/current readings.
$ws2.Range("A1").Value = "Throttle"
$ws2.Range("B1").Value = "Thrust"
$ws2.Range("C1").Value = "Current"

# Match the saved selection/active-sheet state.
[void]$ws1.Range("G10").Select()
[void]$ws2.Range("B4").Select()
[void]$ws2.Activate()
